{"js": "// Update each two-digit \u00f7 one-digit division fact table cell to the new\n// value. The mapping below is derived 1:1 from the authoritative OOXML\n// diff: every \"old\" expression string occurs exactly once in the\n// document, so an exact, case-sensitive whole-text search+replace is\n// unambiguous and independent of cell/row ordering.\nconst replacements = [\n  [\"79\u00f75=15, 4\", \"68\u00f72=34, 0\"],\n  [\"50\u00f74=12, 2\", \"81\u00f78=10, 1\"],\n  [\"54\u00f73=18, 0\", \"64\u00f73=21, 1\"],\n  [\"40\u00f76=6, 4\", \"49\u00f72=24, 1\"],\n  [\"54\u00f77=7, 5\", \"57\u00f77=8, 1\"],\n  [\"47\u00f77=6, 5\", \"35\u00f78=4, 3\"],\n  [\"32\u00f73=10, 2\", \"91\u00f74=22, 3\"],\n  [\"38\u00f74=9, 2\", \"11\u00f79=1, 2\"],\n  [\"27\u00f76=4, 3\", \"41\u00f79=4, 5\"],\n  [\"75\u00f72=37, 1\", \"97\u00f78=12, 1\"],\n  [\"59\u00f78=7, 3\", \"57\u00f74=14, 1\"],\n  [\"26\u00f77=3, 5\", \"77\u00f75=15, 2\"],\n  [\"47\u00f75=9, 2\", \"42\u00f76=7, 0\"],\n  [\"66\u00f75=13, 1\", \"99\u00f73=33, 0\"],\n  [\"74\u00f79=8, 2\", \"63\u00f79=7, 0\"],\n  [\"91\u00f79=10, 1\", \"73\u00f79=8, 1\"],\n  [\"46\u00f78=5, 6\", \"67\u00f74=16, 3\"],\n  [\"89\u00f77=12, 5\", \"45\u00f76=7, 3\"],\n  [\"28\u00f73=9, 1\", \"38\u00f73=12, 2\"],\n  [\"44\u00f75=8, 4\", \"42\u00f79=4, 6\"],\n  [\"41\u00f78=5, 1\", \"30\u00f73=10, 0\"],\n  [\"38\u00f76=6, 2\", \"60\u00f74=15, 0\"],\n  [\"58\u00f76=9, 4\", \"85\u00f79=9, 4\"],\n  [\"51\u00f77=7, 2\", \"70\u00f78=8, 6\"],\n  [\"78\u00f76=13, 0\", \"64\u00f79=7, 1\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old two-digit-division expression with its new value.\n# The mapping is derived 1:1 from the authoritative OOXML diff -- every\n# \"old\" string is unique in the document, so an exact whole-text\n# find/replace (wdReplaceAll) is unambiguous and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"79\u00f75=15, 4\", \"68\u00f72=34, 0\"),\n    @(\"50\u00f74=12, 2\", \"81\u00f78=10, 1\"),\n    @(\"54\u00f73=18, 0\", \"64\u00f73=21, 1\"),\n    @(\"40\u00f76=6, 4\", \"49\u00f72=24, 1\"),\n    @(\"54\u00f77=7, 5\", \"57\u00f77=8, 1\"),\n    @(\"47\u00f77=6, 5\", \"35\u00f78=4, 3\"),\n    @(\"32\u00f73=10, 2\", \"91\u00f74=22, 3\"),\n    @(\"38\u00f74=9, 2\", \"11\u00f79=1, 2\"),\n    @(\"27\u00f76=4, 3\", \"41\u00f79=4, 5\"),\n    @(\"75\u00f72=37, 1\", \"97\u00f78=12, 1\"),\n    @(\"59\u00f78=7, 3\", \"57\u00f74=14, 1\"),\n    @(\"26\u00f77=3, 5\", \"77\u00f75=15, 2\"),\n    @(\"47\u00f75=9, 2\", \"42\u00f76=7, 0\"),\n    @(\"66\u00f75=13, 1\", \"99\u00f73=33, 0\"),\n    @(\"74\u00f79=8, 2\", \"63\u00f79=7, 0\"),\n    @(\"91\u00f79=10, 1\", \"73\u00f79=8, 1\"),\n    @(\"46\u00f78=5, 6\", \"67\u00f74=16, 3\"),\n    @(\"89\u00f77=12, 5\", \"45\u00f76=7, 3\"),\n    @(\"28\u00f73=9, 1\", \"38\u00f73=12, 2\"),\n    @(\"44\u00f75=8, 4\", \"42\u00f79=4, 6\"),\n    @(\"41\u00f78=5, 1\", \"30\u00f73=10, 0\"),\n    @(\"38\u00f76=6, 2\", \"60\u00f74=15, 0\"),\n    @(\"58\u00f76=9, 4\", \"85\u00f79=9, 4\"),\n    @(\"51\u00f77=7, 2\", \"70\u00f78=8, 6\"),\n    @(\"78\u00f76=13, 0\", \"64\u00f79=7, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.Forward = $true\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
